$wb = $excel.ActiveWorkbook

# --- Fix the "tokyo@admin" hyperlink text on the Login sheet (drop trailing space) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "tokyo@admin"
$ws1.Range("C2").Select() | Out-Null

# --- Add the new "AddSbu" sheet (SBU data) after the existing "Login" sheet ---
$ws2 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets($wb.Worksheets.Count))
$ws2.Name = "AddSbu"

$ws2.Range("A1").Value = $true
$ws2.Range("B1").Value = "CEMENT"
$ws2.Range("C1").Value = "qqqqq"

$ws2.Range("A2").Value = $true
$ws2.Range("B2").Value = "RMC"
$ws2.Range("C2").Value = "qqqqqqqq"

$ws2.Range("A3").Value = $true
$ws2.Range("B3").Value = "Test"
$ws2.Range("C3").Value = "qqqqqqqqq"

$ws2.Activate() | Out-Null
$ws2.Range("A1:C3").Select() | Out-Null
